# The "ppm" header label was missing from D1 (it was blank) and had
# accidentally been placed over column F instead. Column F held per-row
# "ppm" values that don't belong in this report, so the whole column F is
# removed - this naturally shifts sample_size/t_results/significance
# (columns G, H, I) one column to the left, into F, G, H - and the
# correct "ppm" header text is written into the now-corrected D1 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value2 = "ppm"
$ws.Columns("F").Delete()
